$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("J2").Value = 3.8
$ws.Range("X2").Value = 990
$ws.Range("AG2").Value = 1000

# Row 3 changes
$ws.Range("F3").Value = 8.199999999999999
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 1.38
$ws.Range("I3").Value = 1.46
$ws.Range("Q3").Value = 1.6
$ws.Range("Z3").Value = 9.6
$ws.Range("AC3").Value = 13
$ws.Range("AE3").Value = 15
$ws.Range("AG3").Value = 40

# Row 4 changes
$ws.Range("F4").Value = 1.08
$ws.Range("H4").Value = 1.08
$ws.Range("J4").Value = 1.13
